$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '24.096.20'
$ws.Range("E2").Value = '  -3.28%  '
$ws.Range("D3").Value = '1.639.65'
$ws.Range("E3").Value = '  -3.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.65%  '
$ws.Range("E6").Value = '  -0.45%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3876'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -1.74%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3851'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.61%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.000'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '49.73'
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.344'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -6.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08694'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.42%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.62'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -7.23%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.066'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.19%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001284'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.84%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.447'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -5.33%  '
$ws.Range("D17").Value = '1.632.38'
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '95.05'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06894'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '20.72'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.876'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.17%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9996'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.53'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.63%  '
$ws.Range("D24").Value = '24.093.43'
$ws.Range("E24").Value = '  -3.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.330'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.19%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.694'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -5.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '157.37'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.94%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.539'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +6.12%  '
$ws.Range("B30").Value = 'BitcoinCash'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '140.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -5.56%  '
$ws.Range("B31").Value = 'HuobiToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.354'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -11.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.406'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -9.48%  '
$ws.Range("D33").Value = '1.812.18'
$ws.Range("E33").Value = '  -5.22%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.871'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.07989'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.22%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02875'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -7.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2666'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.9448'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -8.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.09174'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.26%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.449'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.877'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -8.28%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7531'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -6.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '13.00'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -6.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.95'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -5.48%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6866'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -5.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.456'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.080'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9990'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08370'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -6.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.260'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -9.18%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '132.22'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.06%  '
